$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Tabelle1") updates ---
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Update the three dates (row 4, row 9, row 10) from 2023-12-26 (45286) to 2024-01-02 (45293)
$ws1.Range("B4").Value = 45293
$ws1.Range("B9").Value = 45293
$ws1.Range("B10").Value = 45293

# Add the new comment on row 10 ("Outputs sind Rechts und haben die Farbe gelb")
$ws1.Range("C10").Value = "Müssen konst blöcke auch?"

# Make Tabelle1 the active sheet/tab again, with E10 selected
$ws1.Activate()
$ws1.Range("E10").Select()
